{"js": "const flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:lastRenderedPageBreak/>\n        <w:t>----------------------------------------------------------------------------------------------</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:jc w:val=\"center\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n          <w:highlight w:val=\"yellow\"/>\n        </w:rPr>\n        <w:t>-VISUALIZAR PROGRAMA</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n          <w:highlight w:val=\"yellow\"/>\n        </w:rPr>\n        <w:t>:</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Estado del CU: Desaprobado (Corregir Errores</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Este CU en realidad no presenta \"errores\", sino que lo marque como Desaprobado para que una vez que contemos con el panel principal del sistema VASPA, se le agregue la funcionalidad a los botones que faltaron.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>-Bot\u00f3n \"Cancelar\"  -&gt; Contiene Errores:</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">Resultado esperado: Regresar a la Pantalla Principal del Sistema </w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Vaspa</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"TableContents\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">Resultado obtenido: Insatisfactorio. </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>No se ha regresado a la pantalla principal del sistema VASPA. El bot\u00f3n no cuenta con ninguna funcionalidad.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">-NOTA: Se debe corregir el nombre de la pesta\u00f1a que se abre del navegador cuando quiero visualizar el programa de una asignatura ya que dicha pesta\u00f1a se abre con el nombre del archivo original como el t\u00edtulo. En mi caso, </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>sub\u00ed</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> el programa de Laboratorio de Programaci\u00f3n (para probar nada m\u00e1s) y la pesta\u00f1a llevaba el nombre \"LP programa 2013\", donde deber\u00eda haber sido por ejemplo: \" prg_1649_016_uarg_pact.pdf\" respetando la carrera y el c\u00f3digo. </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n`;\n\nconst body = context.document.body;\nconst range = body.getRange(\"End\");\nrange.insertOoxml(flatOpcXml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Collapse(0)\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:lastRenderedPageBreak/>\n        <w:t>----------------------------------------------------------------------------------------------</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:jc w:val=\"center\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n          <w:highlight w:val=\"yellow\"/>\n        </w:rPr>\n        <w:t>-VISUALIZAR PROGRAMA</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n          <w:highlight w:val=\"yellow\"/>\n        </w:rPr>\n        <w:t>:</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Estado del CU: Desaprobado (Corregir Errores</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Este CU en realidad no presenta \"errores\", sino que lo marque como Desaprobado para que una vez que contemos con el panel principal del sistema VASPA, se le agregue la funcionalidad a los botones que faltaron.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>-Bot\u00f3n \"Cancelar\"  -&gt; Contiene Errores:</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">Resultado esperado: Regresar a la Pantalla Principal del Sistema </w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Vaspa</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"TableContents\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">Resultado obtenido: Insatisfactorio. </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>No se ha regresado a la pantalla principal del sistema VASPA. El bot\u00f3n no cuenta con ninguna funcionalidad.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">-NOTA: Se debe corregir el nombre de la pesta\u00f1a que se abre del navegador cuando quiero visualizar el programa de una asignatura ya que dicha pesta\u00f1a se abre con el nombre del archivo original como el t\u00edtulo. En mi caso, </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>sub\u00ed</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> el programa de Laboratorio de Programaci\u00f3n (para probar nada m\u00e1s) y la pesta\u00f1a llevaba el nombre \"LP programa 2013\", donde deber\u00eda haber sido por ejemplo: \" prg_1649_016_uarg_pact.pdf\" respetando la carrera y el c\u00f3digo. </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:b/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n$result = $rng.InsertXML($xml)\n"}
